{"js": "// Replace the \"Pascal can change everything in the world and never give up ! \"\n// paragraph with the new Version-Management paragraph, collapsing it down to a\n// single plain run (no leftover run-level formatting such as the eastAsia font\n// hint/lang on the old \"Pascal\" run).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst targetText = \"Pascal can change everything in the world and never give up !\";\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.trim() === targetText) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not find the target paragraph to replace.\");\n}\n\nconst range = target.getRange();\n\n// Clear existing runs/formatting first so the replacement run starts from a\n// clean (no rPr) state instead of inheriting the first run's properties.\nrange.clear();\nawait context.sync();\n\nrange.insertText(\n  \"Version Management also called Version Control or Revision Control, is a means to effectively track and control changes to a collection of related entities.\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n", "ps1": "# Replace the \"Pascal can change everything in the world and never give up ! \"\n# paragraph with the new Version-Management paragraph, collapsing it down to a\n# single plain run (no leftover run-level formatting such as the eastAsia font\n# hint/lang carried by the old \"Pascal\" run).\n\n$d = $word.ActiveDocument\n$targetText = \"Pascal can change everything in the world and never give up !\"\n$newText = \"Version Management also called Version Control or Revision Control, is a means to effectively track and control changes to a collection of related entities.\"\n\nforeach ($p in $d.Paragraphs) {\n    $pRange = $p.Range\n    if ($pRange.Text.Trim() -eq $targetText) {\n        # Range over just the paragraph's text, excluding the trailing\n        # paragraph mark, so we don't disturb the paragraph itself.\n        $full = $d.Range($pRange.Start, $pRange.End - 1)\n        # Delete first so the replacement text isn't typed using the\n        # formatting (rPr) inherited from the old first run.\n        $full.Delete()\n        $full.Text = $newText\n        break\n    }\n}\n"}
